# backlog_logiciel_V1.xlsx - "mise a jour doc"
#
# 1. Fix a typo in the backlog: "défi" -> "défis" (missing plural "s")
# 2. Center-align the data table (columns B:H across rows 1:8)
# 3. Re-select cell A2 (matches the saved selection state) and set the zoom
#    level used while reviewing the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Correct the wording of the first backlog item ------------------
$ws.Range("A2").Value = "déposer les vidéeos de défis"

# --- 2. Center every value in the table (creates/applies the new style) -
$ws.Range("B1:H8").HorizontalAlignment = -4108   # xlCenter

# --- 3. Restore view state: selection + zoom ----------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.Zoom = 130
